$d = $word.ActiveDocument

$replacements = @(
    ,@("2025-06-16 Monday", "2025-06-17 Tuesday")
    ,@("13-9=4", "76-40=36")
    ,@("47-22=25", "51-40=11")
    ,@("52-50=2", "45+4=49")
    ,@("85-66=19", "5+9=14")
    ,@("77-12=65", "62-23=39")
    ,@("45+53=98", "66-38=28")
    ,@("74-52=22", "63+2=65")
    ,@("15+79=94", "20-0=20")
    ,@("48-46=2", "91-17=74")
    ,@("72-65=7", "59-24=35")
    ,@("56-3=53", "74-26=48")
    ,@("85-4=81", "23+38=61")
    ,@("17+63=80", "62+20=82")
    ,@("63+3=66", "13+69=82")
    ,@("37-23=14", "24+20=44")
    ,@("80-62=18", "61-7=54")
    ,@("17+78=95", "90-18=72")
    ,@("16+57=73", "57-31=26")
    ,@("0+20=20", "50-38=12")
    ,@("25+31=56", "63+35=98")
    ,@("29+10=39", "39+48=87")
    ,@("17-16=1", "67-1=66")
    ,@("54-21=33", "2+35=37")
    ,@("53-49=4", "38+3=41")
    ,@("33-2=31", "10+78=88")
    ,@("51+36=87", "37+11=48")
    ,@("62-21=41", "89-3=86")
    ,@("10+67=77", "95-83=12")
    ,@("6+84=90", "80-12=68")
    ,@("52-16=36", "49-1=48")
    ,@("17+49=66", "39-29=10")
    ,@("3+61=64", "97-76=21")
    ,@("2+8=10", "33+48=81")
    ,@("51+10=61", "30-5=25")
    ,@("90-62=28", "40+16=56")
    ,@("0+38=38", "88-30=58")
    ,@("30+6=36", "41-13=28")
    ,@("75-31=44", "46-36=10")
    ,@("24+36=60", "71-3=68")
    ,@("21+37=58", "0+27=27")
    ,@("41+9=50", "54+38=92")
    ,@("52-4=48", "71-15=56")
    ,@("49+42=91", "92+3=95")
    ,@("27+64=91", "20+73=93")
    ,@("37+6=43", "78+10=88")
    ,@("58+1=59", "58+38=96")
    ,@("28+32=60", "10+70=80")
    ,@("95-87=8", "29+53=82")
    ,@("23+23=46", "38-11=27")
    ,@("68-53=15", "54+26=80")
    ,@("86-81=5", "26+38=64")
    ,@("51+13=64", "45-7=38")
    ,@("61+10=71", "50-34=16")
    ,@("96-14=82", "94-65=29")
    ,@("84-69=15", "28+62=90")
    ,@("61+7=68", "44-27=17")
    ,@("80-34=46", "8+40=48")
    ,@("68-10=58", "40+14=54")
    ,@("11-5=6", "97-65=32")
    ,@("78-51=27", "99-85=14")
    ,@("73-58=15", "3+35=38")
    ,@("20-11=9", "46-33=13")
    ,@("74-39=35", "84-20=64")
    ,@("47+16=63", "50-8=42")
    ,@("94-45=49", "75+9=84")
    ,@("85-42=43", "70+25=95")
    ,@("33+5=38", "21+36=57")
    ,@("63-32=31", "36-29=7")
    ,@("23-23=0", "27-7=20")
    ,@("13+71=84", "34-9=25")
    ,@("63+17=80", "83-26=57")
    ,@("63-20=43", "50+31=81")
    ,@("84-79=5", "52+9=61")
    ,@("13+14=27", "10+64=74")
    ,@("66-17=49", "27-3=24")
    ,@("43+39=82", "87-20=67")
    ,@("86+5=91", "66-59=7")
    ,@("4+41=45", "17+12=29")
    ,@("37+5=42", "57-0=57")
    ,@("11-2=9", "36-8=28")
    ,@("66-55=11", "98-47=51")
    ,@("22+3=25", "91-27=64")
    ,@("1+58=59", "89+1=90")
    ,@("88-24=64", "13+59=72")
    ,@("87-85=2", "62+15=77")
    ,@("19+11=30", "58-19=39")
    ,@("14-10=4", "70-8=62")
    ,@("5+13=18", "32+29=61")
    ,@("60-2=58", "49+24=73")
    ,@("45+36=81", "85-17=68")
    ,@("7+59=66", "4+29=33")
    ,@("71-31=40", "34-7=27")
    ,@("26+4=30", "89-62=27")
    ,@("40+1=41", "56-35=21")
    ,@("25+12=37", "44-38=6")
    ,@("38+10=48", "54-2=52")
    ,@("45+31=76", "32+16=48")
    ,@("58+6=64", "94-52=42")
    ,@("67-63=4", "96-67=29")
    ,@("76-29=47", "64+23=87")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done"
